$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 4
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = -7
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 13
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = 4
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = -4
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = -4
$ws.Range("F30").Value = -1
$ws.Range("F31").Value = 5
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = 5
$ws.Range("F35").Value = 6
$ws.Range("F36").Value = -4
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 1
